$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as plain text (it already stores values like
# "37.484.72" / "0.651" as text) so values with trailing zeros or
# thousands separators survive the write unmodified.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.444.68"
$ws.Range("E2").Value = "  +5.67%  "
$ws.Range("D3").Value = "2.052.77"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "252.86"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").Value = "0.650"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("D7").Value = "66.26"
$ws.Range("E7").Value = "  +16.66%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +7.20%  "
$ws.Range("D10").Value = "59.57"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "0.910"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "14.88"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("D15").Value = "2.355.55"
$ws.Range("E15").Value = "  +4.14%  "
$ws.Range("D16").Value = "22.36"
$ws.Range("E16").Value = "  +28.17%  "
$ws.Range("D17").Value = "5.58"
$ws.Range("E17").Value = "  +6.46%  "
$ws.Range("D18").Value = "2.045.87"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("D19").Value = "37.274.62"
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("D20").Value = "73.74"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +4.72%  "
$ws.Range("D22").Value = "5.44"
$ws.Range("E22").Value = "  +6.47%  "
$ws.Range("D23").Value = "240.32"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  +5.82%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  +9.60%  "
$ws.Range("D28").Value = "161.61"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").Value = "20.02"
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("E30").Value = "  +28.29%  "
$ws.Range("D31").Value = "5.31"
$ws.Range("E31").Value = "  +9.78%  "
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  +9.52%  "
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  +10.08%  "
$ws.Range("D35").Value = "0.0625"
$ws.Range("E35").Value = "  +6.43%  "
$ws.Range("D36").Value = "2.47"
$ws.Range("E36").Value = "  +5.07%  "
$ws.Range("E37").Value = "  +5.29%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "6.01"
$ws.Range("E39").Value = "  +17.19%  "
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  +34.00%  "
$ws.Range("D41").Value = "0.105"
$ws.Range("E41").Value = "  +18.62%  "
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("D47").Value = "96.54"
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("D48").Value = "8.00"
$ws.Range("E48").Value = "  +7.25%  "
$ws.Range("D49").Value = "1.425.13"
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "46.89"
$ws.Range("E51").Value = "  +1.93%  "

# Row 45 <-> 46 swap (VeChain / InjectiveProtocol)
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0219"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "17.19"
$ws.Range("E46").Value = "  +9.22%  "
